# Generate Report for Handoff
#
# For the "b" file (row 3 in each sheet) mark the item as "Ready for
# handoff" and record the newly produced handoff package for both the
# zh-cn and de-de locales.

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

# ---------------------------------------------------------------------
# Overview sheet: update the zh-cn / de-de status columns for b.md
# ---------------------------------------------------------------------
$wsOverview.Range("B3").Value = "Ready for handoff"
$wsOverview.Range("C3").Value = "Ready for handoff"

# ---------------------------------------------------------------------
# Helper to update a locale detail sheet (zh-cn / de-de) for the row
# belonging to b.md (row 3): status, handoff file name/link, handoff
# datetime.
# ---------------------------------------------------------------------
function Update-LocaleSheet($ws, $handoffFile, $handoffDate) {
    $ws.Range("B3").Value = "Ready for handoff"
    $ws.Range("C3").Value = $handoffFile
    $ws.Range("D3").Value = $handoffDate

    $hyperlink = $null
    foreach ($hl in $ws.Hyperlinks) {
        if ($hl.Range.Address() -eq '$C$3') {
            $hyperlink = $hl
        }
    }
    if ($hyperlink -ne $null) {
        $hyperlink.TextToDisplay = $handoffFile
    }
}

Update-LocaleSheet $wsZhCn "b.63290e5768f688058c7b37413b0a5c26c308f864.zh-cn.xlf" "2016-03-10 08:53:07"
Update-LocaleSheet $wsDeDe "b.63290e5768f688058c7b37413b0a5c26c308f864.de-de.xlf" "2016-03-10 08:53:14"
